$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices and volume %) per row.
# Column D (Price) cells are numeric-looking text; force text format so
# exact string representation (e.g. trailing zeros, thousand-separators) is preserved.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.199.65'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.250.08'
$ws.Range("E3").Value = '  +5.39%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.55'
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.25'
$ws.Range("E6").Value = '  +7.13%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.239.08'
$ws.Range("E8").Value = '  +5.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +3.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.09'
$ws.Range("E10").Value = '  +10.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.165'
$ws.Range("E11").Value = '  +4.01%  '
$ws.Range("E12").Value = '  +4.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.87'
$ws.Range("E13").Value = '  +3.90%  '
$ws.Range("E14").Value = '  +4.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.755.86'
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '558.89'
$ws.Range("E16").Value = '  +12.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.249.51'
$ws.Range("E17").Value = '  +2.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.247.52'
$ws.Range("E18").Value = '  +5.26%  '
$ws.Range("E19").Value = '  +2.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.13'
$ws.Range("E20").Value = '  +5.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.46'
$ws.Range("E21").Value = '  +3.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.746'
$ws.Range("E22").Value = '  +7.11%  '
$ws.Range("E23").Value = '  +8.41%  '
$ws.Range("E24").Value = '  +6.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.21'
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.42'
$ws.Range("E27").Value = '  +18.07%  '
$ws.Range("E28").Value = '  +5.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +6.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.92'
$ws.Range("E30").Value = '  +5.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.78'
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("E33").Value = '  +4.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '572.53'
$ws.Range("E34").Value = '  +10.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.80'
$ws.Range("E35").Value = '  +4.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.42'
$ws.Range("E36").Value = '  +5.48%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.48'
$ws.Range("E37").Value = '  +2.92%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0457'
$ws.Range("E38").Value = '  +12.06%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0868'
$ws.Range("E39").Value = '  +6.93%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.06'
$ws.Range("E40").Value = '  +12.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.128'
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.157.03'
$ws.Range("E42").Value = '  +6.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.64'
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.275'
$ws.Range("E44").Value = '  +9.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("E45").Value = '  +5.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.68'
$ws.Range("E46").Value = '  +4.05%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₃0558'
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.33'
$ws.Range("E49").Value = '  +3.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.113'
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.27'
$ws.Range("E51").Value = '  +7.75%  '
